# edit.ps1
# Applies two changes described by the source diff:
#   1. Slide 1 subtitle: "Dirk Riehle, Univ. Erlangen" -> "Dirk Riehle, FAU Erlangen"
#   2. Theme colour scheme used by the deck (ppt/theme/theme1.xml, reached via
#      the slide master's Theme object) switches from the "Simple Light" /
#      "ADAP Slides Template" palette to the "Default" palette.

$p = $ppt.ActivePresentation

# --- 1. Fix author affiliation on the title slide -------------------------
$s1 = $p.Slides.Item(1)
$subtitleShape = $s1.Shapes.Item(2)
$run = $subtitleShape.TextFrame.TextRange.Paragraphs(1).Runs(1)
$run.Text = "Dirk Riehle, FAU Erlangen"

# --- 2. Swap the theme colour palette --------------------------------------
function ConvertTo-BgrColor([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return ($b -shl 16) -bor ($g -shl 8) -bor $r
}

# Theme colour order exposed by ThemeColorScheme.Colors(i):
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3,
#   8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink
$defaultPalette = @(
    "000000",
    "FFFFFF",
    "158158",
    "F3F3F3",
    "058DC7",
    "50B432",
    "ED561B",
    "EDEF00",
    "24CBE5",
    "64E572",
    "2200CC",
    "551A8B"
)

$themeColors = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Colors($i).RGB = ConvertTo-BgrColor $defaultPalette[$i - 1]
}
